$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "is_active" column (F) was stored as a TRUE() formula evaluating to
# the number 1. Replace it with the literal text string "TRUE" (shared
# string), keeping the existing Text-formatted style (s="1") on the cells.
$rng = $ws.Range("F2:F11")
$rng.Formula = '="TRUE"'
$rng.Copy()
$rng.PasteSpecial(-4163)   # xlPasteValues - collapse formula to a static value

# Move the active selection to H10, as recorded in the saved workbook view.
[void]$ws.Range("H10").Select()
